$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = 39583
$ws.Cells.Item(2, 2).Value2 = 2008
$ws.Cells.Item(2, 3).ClearContents() | Out-Null
$ws.Cells.Item(2, 4).Value2 = 2009
$ws.Cells.Item(2, 5).Value2 = -0.3422723562191532

$ws.Cells.Item(3, 1).Value2 = 39765
$ws.Cells.Item(3, 2).Value2 = 2008
$ws.Cells.Item(3, 3).ClearContents() | Out-Null
$ws.Cells.Item(3, 4).Value2 = 2009
$ws.Cells.Item(3, 5).Value2 = 0.4944284391569687

$ws.Cells.Item(4, 1).Value2 = 39948
$ws.Cells.Item(4, 2).Value2 = 2009
$ws.Cells.Item(4, 3).Value2 = -0.5037688924316441
$ws.Cells.Item(4, 4).Value2 = 2010
$ws.Cells.Item(4, 5).Value2 = -0.467076459743887

$ws.Cells.Item(5, 1).Value2 = 40130
$ws.Cells.Item(5, 2).Value2 = 2009
$ws.Cells.Item(5, 3).Value2 = -0.5555135891318952
$ws.Cells.Item(5, 4).Value2 = 2010
$ws.Cells.Item(5, 5).Value2 = -0.4782015746048418

$ws.Cells.Item(6, 1).Value2 = 40310
$ws.Cells.Item(6, 2).Value2 = 2010
$ws.Cells.Item(6, 3).Value2 = 0.169534172659791
$ws.Cells.Item(6, 4).Value2 = 2011
$ws.Cells.Item(6, 5).Value2 = -0.03968684591929561

$ws.Cells.Item(7, 1).Value2 = 40494
$ws.Cells.Item(7, 2).Value2 = 2010
$ws.Cells.Item(7, 3).Value2 = 0.8442071301477228
$ws.Cells.Item(7, 4).Value2 = 2011
$ws.Cells.Item(7, 5).Value2 = 1.324233212457782

$ws.Cells.Item(8, 1).Value2 = 40676
$ws.Cells.Item(8, 2).Value2 = 2011
$ws.Cells.Item(8, 3).Value2 = 1.032338390744236
$ws.Cells.Item(8, 4).Value2 = 2012
$ws.Cells.Item(8, 5).Value2 = 0.3765075513336269

$ws.Cells.Item(9, 1).Value2 = 40862
$ws.Cells.Item(9, 2).Value2 = 2011
$ws.Cells.Item(9, 3).Value2 = 1.122475521884692
$ws.Cells.Item(9, 4).Value2 = 2012
$ws.Cells.Item(9, 5).Value2 = 0.7478380109886329

$ws.Cells.Item(10, 1).Value2 = 41044
$ws.Cells.Item(10, 2).Value2 = 2012
$ws.Cells.Item(10, 3).Value2 = -0.1937612543835177
$ws.Cells.Item(10, 4).Value2 = 2013
$ws.Cells.Item(10, 5).Value2 = 0.6066710853121382

$ws.Cells.Item(11, 1).Value2 = 41228
$ws.Cells.Item(11, 2).Value2 = 2012
$ws.Cells.Item(11, 3).Value2 = -0.578174579726376
$ws.Cells.Item(11, 4).Value2 = 2013
$ws.Cells.Item(11, 5).Value2 = -0.2445716668737163

$ws.Cells.Item(12, 1).Value2 = 41409
$ws.Cells.Item(12, 2).Value2 = 2013
$ws.Cells.Item(12, 3).Value2 = -0.8331679621937482
$ws.Cells.Item(12, 4).Value2 = 2014
$ws.Cells.Item(12, 5).Value2 = -0.3970496740026364

$ws.Cells.Item(13, 1).Value2 = 41592
$ws.Cells.Item(13, 2).Value2 = 2013
$ws.Cells.Item(13, 3).Value2 = -0.7492845378401558
$ws.Cells.Item(13, 4).Value2 = 2014
$ws.Cells.Item(13, 5).Value2 = -0.2617076051026235

$ws.Cells.Item(14, 1).Value2 = 41774
$ws.Cells.Item(14, 2).Value2 = 2014
$ws.Cells.Item(14, 3).Value2 = 0.1273541662098365
$ws.Cells.Item(14, 4).Value2 = 2015
$ws.Cells.Item(14, 5).Value2 = -0.4865818826308876

$ws.Cells.Item(15, 1).Value2 = 41957
$ws.Cells.Item(15, 2).Value2 = 2014
$ws.Cells.Item(15, 3).Value2 = 0.2751437421933511
$ws.Cells.Item(15, 4).Value2 = 2015
$ws.Cells.Item(15, 5).Value2 = -0.100009932057743

$ws.Cells.Item(16, 1).Value2 = 42137
$ws.Cells.Item(16, 2).Value2 = 2015
$ws.Cells.Item(16, 3).Value2 = -0.08273351073040391
$ws.Cells.Item(16, 4).Value2 = 2016
$ws.Cells.Item(16, 5).Value2 = -0.01252079199893785

$ws.Cells.Item(17, 1).Value2 = 42321
$ws.Cells.Item(17, 2).Value2 = 2015
$ws.Cells.Item(17, 3).Value2 = 0.07468705617190707
$ws.Cells.Item(17, 4).Value2 = 2016
$ws.Cells.Item(17, 5).Value2 = 0.3000376062062493

$ws.Cells.Item(18, 1).Value2 = 42503
$ws.Cells.Item(18, 2).Value2 = 2016
$ws.Cells.Item(18, 3).Value2 = -0.1151820594382569
$ws.Cells.Item(18, 4).Value2 = 2017
$ws.Cells.Item(18, 5).Value2 = 0.07348980370169844

$ws.Cells.Item(19, 1).Value2 = 42689
$ws.Cells.Item(19, 2).Value2 = 2016
$ws.Cells.Item(19, 3).Value2 = -0.05493014849097255
$ws.Cells.Item(19, 4).Value2 = 2017
$ws.Cells.Item(19, 5).Value2 = 0.1740313431290996

$ws.Cells.Item(20, 1).Value2 = 42867
$ws.Cells.Item(20, 2).Value2 = 2017
$ws.Cells.Item(20, 3).Value2 = 0.157394256377752
$ws.Cells.Item(20, 4).Value2 = 2018
$ws.Cells.Item(20, 5).Value2 = 0.06194937150048041

$ws.Cells.Item(21, 1).Value2 = 43053
$ws.Cells.Item(21, 2).Value2 = 2017
$ws.Cells.Item(21, 3).Value2 = 0.2820931576894115
$ws.Cells.Item(21, 4).Value2 = 2018
$ws.Cells.Item(21, 5).Value2 = 0.3390041783450259

$ws.Cells.Item(22, 1).Value2 = 43145
$ws.Cells.Item(22, 2).Value2 = 2018
$ws.Cells.Item(22, 3).Value2 = 0.2895071529679827
$ws.Cells.Item(22, 4).Value2 = 2019
$ws.Cells.Item(22, 5).Value2 = 0.1190486358061627

$ws.Cells.Item(23, 1).Value2 = 43235
$ws.Cells.Item(23, 2).Value2 = 2018
$ws.Cells.Item(23, 3).Value2 = 0.3426151435189873
$ws.Cells.Item(23, 4).Value2 = 2019
$ws.Cells.Item(23, 5).Value2 = 0.1686730364466316

$ws.Cells.Item(24, 1).Value2 = 43326
$ws.Cells.Item(24, 2).Value2 = 2018
$ws.Cells.Item(24, 3).Value2 = 0.187607693984293
$ws.Cells.Item(24, 4).Value2 = 2019
$ws.Cells.Item(24, 5).Value2 = 0.02358181985058216

$ws.Cells.Item(25, 1).Value2 = 43418
$ws.Cells.Item(25, 2).Value2 = 2018
$ws.Cells.Item(25, 3).Value2 = 0.2343541283920114
$ws.Cells.Item(25, 4).Value2 = 2019
$ws.Cells.Item(25, 5).Value2 = 0.2210188332817387

$ws.Cells.Item(26, 1).Value2 = 43510
$ws.Cells.Item(26, 2).Value2 = 2019
$ws.Cells.Item(26, 3).Value2 = 0.3712693419885671
$ws.Cells.Item(26, 4).Value2 = 2020
$ws.Cells.Item(26, 5).Value2 = 0.3003887663412641

$ws.Cells.Item(27, 1).Value2 = 43600
$ws.Cells.Item(27, 2).Value2 = 2019
$ws.Cells.Item(27, 3).Value2 = 0.1598952850611068
$ws.Cells.Item(27, 4).Value2 = 2020
$ws.Cells.Item(27, 5).Value2 = 0.132465972367557

$ws.Cells.Item(28, 1).Value2 = 43691
$ws.Cells.Item(28, 2).Value2 = 2019
$ws.Cells.Item(28, 3).Value2 = -0.02097628618118463
$ws.Cells.Item(28, 4).Value2 = 2020
$ws.Cells.Item(28, 5).Value2 = -0.1339126064348495

$ws.Cells.Item(29, 1).Value2 = 43783
$ws.Cells.Item(29, 2).Value2 = 2019
$ws.Cells.Item(29, 3).Value2 = -0.009430310228020211
$ws.Cells.Item(29, 4).Value2 = 2020
$ws.Cells.Item(29, 5).Value2 = -0.09571633453315798

$ws.Cells.Item(30, 1).Value2 = 43875
$ws.Cells.Item(30, 2).Value2 = 2020
$ws.Cells.Item(30, 3).Value2 = -0.03054415496863694
$ws.Cells.Item(30, 4).Value2 = 2021
$ws.Cells.Item(30, 5).Value2 = 0.035598638033707

$ws.Cells.Item(31, 1).Value2 = 43966
$ws.Cells.Item(31, 2).Value2 = 2020
$ws.Cells.Item(31, 3).Value2 = -0.4923796969465988
$ws.Cells.Item(31, 4).Value2 = 2021
$ws.Cells.Item(31, 5).Value2 = -0.3459257698102514

$ws.Cells.Item(32, 1).Value2 = 44068
$ws.Cells.Item(32, 2).Value2 = 2020
$ws.Cells.Item(32, 3).Value2 = -2.657403949513992
$ws.Cells.Item(32, 4).Value2 = 2021
$ws.Cells.Item(32, 5).Value2 = -2.423328265806446

$ws.Cells.Item(33, 1).Value2 = 44159
$ws.Cells.Item(33, 2).Value2 = 2020
$ws.Cells.Item(33, 3).Value2 = -2.657403949513992
$ws.Cells.Item(33, 4).Value2 = 2021
$ws.Cells.Item(33, 5).Value2 = -1.49562970548649

$ws.Cells.Item(34, 1).Value2 = 44251
$ws.Cells.Item(34, 2).Value2 = 2021
$ws.Cells.Item(34, 3).Value2 = -0.2885033948250459
$ws.Cells.Item(34, 4).Value2 = 2022
$ws.Cells.Item(34, 5).Value2 = -0.599348850912329

$ws.Cells.Item(35, 1).Value2 = 44341
$ws.Cells.Item(35, 2).Value2 = 2021
$ws.Cells.Item(35, 3).Value2 = -0.5121403324772844
$ws.Cells.Item(35, 4).Value2 = 2022
$ws.Cells.Item(35, 5).Value2 = -1.073589070820447

$ws.Cells.Item(36, 1).Value2 = 44432
$ws.Cells.Item(36, 2).Value2 = 2021
$ws.Cells.Item(36, 3).Value2 = -0.3096364143617802
$ws.Cells.Item(36, 4).Value2 = 2022
$ws.Cells.Item(36, 5).Value2 = -0.3018961902350958

$ws.Cells.Item(37, 1).Value2 = 44525
$ws.Cells.Item(37, 2).Value2 = 2021
$ws.Cells.Item(37, 3).Value2 = -0.3096364143617802
$ws.Cells.Item(37, 4).Value2 = 2022
$ws.Cells.Item(37, 5).Value2 = -0.1048501255800471

$ws.Cells.Item(38, 1).Value2 = 44617
$ws.Cells.Item(38, 2).Value2 = 2022
$ws.Cells.Item(38, 3).Value2 = -0.2588455356339781
$ws.Cells.Item(38, 4).Value2 = 2023
$ws.Cells.Item(38, 5).Value2 = -1.151120647939763

$ws.Cells.Item(39, 1).Value2 = 44706
$ws.Cells.Item(39, 2).Value2 = 2022
$ws.Cells.Item(39, 3).Value2 = -0.2454721753057276
$ws.Cells.Item(39, 4).Value2 = 2023
$ws.Cells.Item(39, 5).Value2 = -1.238905350026021

$ws.Cells.Item(40, 1).Value2 = 44798
$ws.Cells.Item(40, 2).Value2 = 2022
$ws.Cells.Item(40, 3).Value2 = -0.1730430455425092
$ws.Cells.Item(40, 4).Value2 = 2023
$ws.Cells.Item(40, 5).Value2 = -0.8943276391025989

$ws.Cells.Item(41, 1).Value2 = 44890
$ws.Cells.Item(41, 2).Value2 = 2022
$ws.Cells.Item(41, 3).Value2 = -0.1730430455425092
$ws.Cells.Item(41, 4).Value2 = 2023
$ws.Cells.Item(41, 5).Value2 = 0.9692952624595019

$ws.Cells.Item(42, 1).Value2 = 44981
$ws.Cells.Item(42, 2).Value2 = 2023
$ws.Cells.Item(42, 3).Value2 = 0.8644693227634503
$ws.Cells.Item(42, 4).Value2 = 2024
$ws.Cells.Item(42, 5).Value2 = 0.7652102000489602

$ws.Cells.Item(43, 1).Value2 = 45071
$ws.Cells.Item(43, 2).Value2 = 2023
$ws.Cells.Item(43, 3).Value2 = 0.7038634017465073
$ws.Cells.Item(43, 4).Value2 = 2024
$ws.Cells.Item(43, 5).Value2 = 0.4406734233171727

$ws.Cells.Item(44, 1).Value2 = 45163
$ws.Cells.Item(44, 2).Value2 = 2023
$ws.Cells.Item(44, 3).Value2 = 0.6376744206510576
$ws.Cells.Item(44, 4).Value2 = 2024
$ws.Cells.Item(44, 5).Value2 = -0.09128981027868299

$ws.Cells.Item(45, 1).Value2 = 45254
$ws.Cells.Item(45, 2).Value2 = 2023
$ws.Cells.Item(45, 3).Value2 = 0.6376744206510576
$ws.Cells.Item(45, 4).Value2 = 2024
$ws.Cells.Item(45, 5).Value2 = 0.1544084105021826

$ws.Cells.Item(46, 1).Value2 = 45345
$ws.Cells.Item(46, 2).Value2 = 2024
$ws.Cells.Item(46, 3).Value2 = 0.001611361207976003
$ws.Cells.Item(46, 4).Value2 = 2025
$ws.Cells.Item(46, 5).Value2 = 0.01799217181808199

$ws.Cells.Item(47, 1).Value2 = 45436
$ws.Cells.Item(47, 2).Value2 = 2024
$ws.Cells.Item(47, 3).Value2 = 0.20168190406884
$ws.Cells.Item(47, 4).Value2 = 2025
$ws.Cells.Item(47, 5).Value2 = 0.3613321345859122

$ws.Cells.Item(48, 1).Value2 = 45534
$ws.Cells.Item(48, 2).Value2 = 2024
$ws.Cells.Item(48, 3).Value2 = 0.1856341247700399
$ws.Cells.Item(48, 4).Value2 = 2025
$ws.Cells.Item(48, 5).Value2 = 0.4108497965175983

$ws.Cells.Item(49, 1).Value2 = 45618
$ws.Cells.Item(49, 2).Value2 = 2024
$ws.Cells.Item(49, 3).Value2 = 0.1856341247700399
$ws.Cells.Item(49, 4).Value2 = 2025
$ws.Cells.Item(49, 5).Value2 = 0.3997355152047577

$ws.Cells.Item(50, 1).Value2 = 45713
$ws.Cells.Item(50, 2).Value2 = 2025
$ws.Cells.Item(50, 3).Value2 = 0.08117592553187336
$ws.Cells.Item(50, 4).Value2 = 2026
$ws.Cells.Item(50, 5).Value2 = 0.3069836986764551

$ws.Cells.Item(51, 1).Value2 = 45800
$ws.Cells.Item(51, 2).Value2 = 2025
$ws.Cells.Item(51, 3).Value2 = -0.06418790329880686
$ws.Cells.Item(51, 4).Value2 = 2026
$ws.Cells.Item(51, 5).Value2 = -0.08988642825158433

$ws.Cells.Item(52, 1).Value2 = 45891
$ws.Cells.Item(52, 2).Value2 = 2025
$ws.Cells.Item(52, 3).Value2 = -0.09450306168263811
$ws.Cells.Item(52, 4).Value2 = 2026
$ws.Cells.Item(52, 5).Value2 = -0.3000102673190841

$ws.Rows.Item(53).Delete() | Out-Null
